$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Linear)
$ws.Range("B2").Value = 1.934404800538057
$ws.Range("C2").Value = 3.386296686189082
$ws.Range("D2").Value = 0.6975977593485854

# Row 3 (Decision Tree)
$ws.Range("B3").Value = 2.241367533100733
$ws.Range("C3").Value = 4.125
$ws.Range("D3").Value = 0.4549355386378225

# Row 4 (Random Forest)
$ws.Range("B4").Value = 1.571478190491265
$ws.Range("C4").Value = 2.142957746478873
$ws.Range("D4").Value = 0.8682868937538168

# Row 5 (Lasso)
$ws.Range("B5").Value = 1.920857547837128
$ws.Range("C5").Value = 3.306988653938654
$ws.Range("D5").Value = 0.7059804602205166
